$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 8914.182000000001
$ws.Range("J112").Value = 8914.182000000001
$ws.Range("L112").Value = 26742.546
$ws.Range("N112").Value = -28958.546
$ws.Range("H115").Value = 3497.1667
$ws.Range("I115").Value = 3906.6
$ws.Range("K115").Value = 11719.8
$ws.Range("M115").Value = -10152.8
$ws.Range("H138").Value = 4572.2573
$ws.Range("I138").Value = 3432.9092
$ws.Range("K138").Value = 10298.7276
$ws.Range("M138").Value = -5158.7276

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3950.4348
$ws.Range("I61").Value = 3567.2856
$ws.Range("K61").Value = 3567.2856
$ws.Range("M61").Value = -3355.2856
$ws.Range("H74").Value = 8372.459999999999
$ws.Range("I74").Value = 6155.838
$ws.Range("K74").Value = 6155.838
$ws.Range("M74").Value = -5281.838
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H77").Value = 8372.459999999999
$ws.Range("I77").Value = 6155.838
$ws.Range("K77").Value = 30779.19
$ws.Range("M77").Value = -26411.19
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H122").Value = 5477.354
$ws.Range("I122").Value = 5540.5854
$ws.Range("J122").Value = 5107
$ws.Range("K122").Value = 16621.7562
$ws.Range("L122").Value = 15321
$ws.Range("M122").Value = -14171.7562
$ws.Range("N122").Value = -20221
$ws.Range("H136").Value = 3950.4348
$ws.Range("I136").Value = 3567.2856
$ws.Range("K136").Value = 10701.8568
$ws.Range("M136").Value = -8151.856800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8661.823
$ws.Range("I86").Value = 6738.0835
$ws.Range("K86").Value = 6738.0835
$ws.Range("M86").Value = -5615.0835
$ws.Range("H88").Value = 17086.666
$ws.Range("J88").Value = 14879.6
$ws.Range("L88").Value = 14879.6
$ws.Range("N88").Value = -15691.6
$ws.Range("H89").Value = 8661.823
$ws.Range("I89").Value = 6738.0835
$ws.Range("K89").Value = 33690.4175
$ws.Range("M89").Value = -28074.4175
$ws.Range("H91").Value = 17086.666
$ws.Range("J91").Value = 14879.6
$ws.Range("L91").Value = 14879.6
$ws.Range("N91").Value = -17687.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 548000.9399999999
$ws.Range("I31").Value = 9220.1875
$ws.Range("K31").Value = 9220.1875
$ws.Range("M31").Value = -8925.1875
$ws.Range("H34").Value = 548000.9399999999
$ws.Range("I34").Value = 9220.1875
$ws.Range("K34").Value = 9220.1875
$ws.Range("M34").Value = -9018.1875
$ws.Range("H51").Value = 68099.8
$ws.Range("J51").Value = 96166.664
$ws.Range("L51").Value = 96166.664
$ws.Range("N51").Value = -97638.664
$ws.Range("H61").Value = 68099.8
$ws.Range("J61").Value = 96166.664
$ws.Range("L61").Value = 96166.664
$ws.Range("N61").Value = -96862.664
$ws.Range("H62").Value = 2999.75
$ws.Range("I62").Value = 2999.75
$ws.Range("K62").Value = 2999.75
$ws.Range("M62").Value = -2375.75
$ws.Range("H64").Value = 64849.5
$ws.Range("I64").Value = 64849.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 64849.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -64601.5
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 2999.75
$ws.Range("I65").Value = 2999.75
$ws.Range("K65").Value = 14998.75
$ws.Range("M65").Value = -11878.75
$ws.Range("H67").Value = 64849.5
$ws.Range("I67").Value = 64849.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 64849.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -63991.5
$ws.Range("N67").ClearContents()
$ws.Range("H132").Value = 2831.5938
$ws.Range("I132").Value = 2300.923
$ws.Range("J132").Value = 5131.1665
$ws.Range("K132").Value = 6902.768999999999
$ws.Range("L132").Value = 15393.4995
$ws.Range("M132").Value = -4372.768999999999
$ws.Range("N132").Value = -20453.4995
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5500
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 16500
$ws.Range("N64").Value = -17040
$ws.Range("H67").Value = 5500
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 16500
$ws.Range("N67").Value = -18372
$ws.Range("H109").Value = 3299
$ws.Range("I109").Value = 3299
$ws.Range("K109").Value = 9897
$ws.Range("M109").Value = -8857
$ws.Range("H119").Value = 3711
$ws.Range("I119").Value = 3711
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 11133
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -6295
$ws.Range("N119").ClearContents()
$ws.Range("H130").Value = 4524.25
$ws.Range("I130").Value = 4466
$ws.Range("J130").Value = 4699
$ws.Range("K130").Value = 13398
$ws.Range("L130").Value = 14097
$ws.Range("M130").Value = -8378
$ws.Range("N130").Value = -24137
$ws.Range("H131").Value = 2930.087
$ws.Range("J131").Value = 3744.0715
$ws.Range("L131").Value = 11232.2145
$ws.Range("N131").Value = -21312.2145
$ws.Range("H137").Value = 1859.75
$ws.Range("I137").Value = 1646.5
$ws.Range("J137").Value = 2499.5
$ws.Range("K137").Value = 4939.5
$ws.Range("L137").Value = 7498.5
$ws.Range("M137").Value = 160.5
$ws.Range("N137").Value = -17698.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3027.5
$ws.Range("I113").Value = 2713.6667
$ws.Range("K113").Value = 2713.6667
$ws.Range("M113").Value = -543.6667000000002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H136").Value = 5017.645
$ws.Range("I136").Value = 5283.6665
$ws.Range("K136").Value = 15850.9995
$ws.Range("M136").Value = -13300.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws.Range("H136").Value = 8398.323
$ws.Range("I136").Value = 12414.728
$ws.Range("J136").Value = 6477.4346
$ws.Range("K136").Value = 37244.18399999999
$ws.Range("L136").Value = 19432.3038
$ws.Range("M136").Value = -34694.18399999999
$ws.Range("N136").Value = -24532.3038
